# Generate Report for Handoff
# Adds a new "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md" row to each of the
# three report sheets (Overview, zh-cn, de-de), growing each table by one row.

$wb = $excel.ActiveWorkbook

$repo      = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/a2a7368f9821e1e503e3f1e00d3e2ae078a153ca/e2e/517ad39b-9fcc-4839-9da2-da4ec69df8e6.md"
$fileName  = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md"
$pathName  = "e2e\517ad39b-9fcc-4839-9da2-da4ec69df8e6.md"
$dateFmt   = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview" -> new row 5
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Cells.Item(5, 1).Value = $fileName                 # A5 File Name
$wsOv.Cells.Item(5, 3).Value = ".md"                      # C5 Extension
$wsOv.Cells.Item(5, 4).Value = ""                         # D5 Publish URL (blank)
$wsOv.Cells.Item(5, 5).Value = "Ready for handoff"        # E5 zh-cn
$wsOv.Cells.Item(5, 6).Value = "Ready for handoff"        # F5 de-de
$wsOv.Cells.Item(5, 7).Value = "2017-02-22 08:26:31"      # G5 Latest HO Xliff Generate Date
$wsOv.Cells.Item(5, 7).NumberFormat = $dateFmt

$wsOv.Hyperlinks.Add($wsOv.Cells.Item(5, 2), $repo, "", "", $pathName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" -> new row 5
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Cells.Item(5, 2).Value  = ".md"                      # B5 File Extension
$wsZh.Cells.Item(5, 3).Value  = "Ready for handoff"        # C5 Status
$wsZh.Cells.Item(5, 4).Value  = "e2e"                      # D5 Source Path
$wsZh.Cells.Item(5, 5).Value  = "ht"                       # E5 Priority
$wsZh.Cells.Item(5, 6).Value  = "'False"                   # F5 Content Duplicate
$wsZh.Cells.Item(5, 7).Value  = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.b354480c3ab030bef42963314777d4b2788f3626.zh-cn.xlf" # G5 Latest Handoff File
$wsZh.Cells.Item(5, 8).Value  = "2017-02-22 08:26:15"      # H5 Latest Handoff Datetime
$wsZh.Cells.Item(5, 8).NumberFormat = $dateFmt
$wsZh.Cells.Item(5, 12).Value = "0001-01-01 00:00:00"      # L5 Latest Handback DateTime
$wsZh.Cells.Item(5, 12).NumberFormat = $dateFmt
$wsZh.Cells.Item(5, 15).Value = "'True"                    # O5 Has metadata
$wsZh.Cells.Item(5, 17).Value = "'False"                   # Q5 Has metadata (duplicate flag col)

$wsZh.Hyperlinks.Add($wsZh.Cells.Item(5, 1), $repo, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" -> new row 5
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Cells.Item(5, 2).Value  = ".md"                      # B5 File Extension
$wsDe.Cells.Item(5, 3).Value  = "Ready for handoff"        # C5 Status
$wsDe.Cells.Item(5, 4).Value  = "e2e"                      # D5 Source Path
$wsDe.Cells.Item(5, 5).Value  = "ht"                       # E5 Priority
$wsDe.Cells.Item(5, 6).Value  = "'False"                   # F5 Content Duplicate
$wsDe.Cells.Item(5, 7).Value  = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.b354480c3ab030bef42963314777d4b2788f3626.de-de.xlf" # G5 Latest Handoff File
$wsDe.Cells.Item(5, 8).Value  = "2017-02-22 08:26:31"      # H5 Latest Handoff Datetime
$wsDe.Cells.Item(5, 8).NumberFormat = $dateFmt
$wsDe.Cells.Item(5, 12).Value = "0001-01-01 00:00:00"      # L5 Latest Handback DateTime
$wsDe.Cells.Item(5, 12).NumberFormat = $dateFmt
$wsDe.Cells.Item(5, 15).Value = "'True"                    # O5 Has metadata
$wsDe.Cells.Item(5, 17).Value = "'False"                   # Q5 Has metadata (duplicate flag col)

$wsDe.Hyperlinks.Add($wsDe.Cells.Item(5, 1), $repo, "", "", $fileName) | Out-Null
